# Daily attendance processing - 2026-02-20 12:08:52 UTC
# Reorders the comma-separated "Recorded By" values in column G so that
# the "2025/2026" token is listed first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = "2025/2026, 2022/2023"
    22 = "2025/2026, 2024/2025"
    23 = "2025/2026, 2023/2024, 2022/2023"
    24 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    27 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    28 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    31 = "2025/2026, 2022/2023"
    50 = "2025/2026, 2024/2025"
    51 = "2025/2026, 2023/2024, 2022/2023"
    52 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    55 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    56 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
